# Update "想去人数" (F column) figures on both the "展览" and "全部类型"
# sheets, which carry duplicate copies of the same exhibition rows.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 691
    $ws.Range("F3").Value = 4012
    $ws.Range("F4").Value = 112
}
